$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) row
$ws.Range("B13").Value = "'34.74"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'14.21"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'48.95"
$ws.Range("D13").Style = "Normal"

# Employment (% of total) row
$ws.Range("B14").Value = "'5.13"
$ws.Range("B14").Style = "Normal"
$ws.Range("D14").Value = "'90.43"
$ws.Range("D14").Style = "Normal"

# Enterprises (% of total) row
$ws.Range("C16").Value = "'29.01"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'99.91"
$ws.Range("D16").Style = "Normal"
